$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit removes the data row for account "005009922" (ANA, balance 3507.85).
# Locate that row dynamically (falls back to the known row 8) and delete it,
# which shifts all subsequent rows up by one.
$target = $ws.Columns.Item(1).Find("005009922")
if ($target) {
    $target.EntireRow.Delete()
} else {
    $ws.Rows.Item(8).Delete()
}
